# Weekly update: a new "Ajo" (garlic) price entry for Terminal Hortofrutícola
# Agro Chillán is inserted as the new row 385 (pushing the existing rows
# 385-399 down to 386-400). The new entry mirrors the prior row 385's
# values except for a newer reporting date (Fecha).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 385; everything below shifts down one
# row (old row 385 -> 386, ..., old row 399 -> 400), and the sheet's used
# range / dimension grows to R400 automatically.
$ws.Rows.Item(385).Insert()

# Populate the newly inserted row 385 with the new weekly record.
$ws.Range("A385").Value = 7
$ws.Range("B385").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C385").Value = "Ñuble"
$ws.Range("D385").Value = 45041
$ws.Range("E385").Value = 16
$ws.Range("F385").Value = 100112003
$ws.Range("G385").Value = "Ajo"
$ws.Range("H385").Value = "Chino"
$ws.Range("I385").Value = "Primera"
$ws.Range("J385").Value = 60
$ws.Range("K385").Value = 17000
$ws.Range("L385").Value = 18000
$ws.Range("M385").Value = 17500
$ws.Range("N385").Value = "$/caja 10 kilos"
$ws.Range("O385").Value = "China"
$ws.Range("P385").Value = 1750
$ws.Range("Q385").Value = 10
$ws.Range("R385").Value = "Hortaliza"
